# Add two new rows of data (ESPESOR 1.2 / "N" and "G") below the existing
# insumos-piezas table, following the same layout as the rows above.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: 1.2 | N | 10013
$ws.Range("A6").Value = 1.2
$ws.Range("B6").Value = "N"
$ws.Range("C6").Value = 10013

# Row 7: 1.2 | G | 20013
$ws.Range("A7").Value = 1.2
$ws.Range("B7").Value = "G"
$ws.Range("C7").Value = 20013

# Match the number formatting/border used by the rest of the table for
# columns A (ESPESOR) and C (code), copied down from the row above.
$ws.Range("A4:A5").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C4:C5").Copy()
$ws.Range("C6:C7").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Leave the selection where the user would have ended up after typing the
# last row (the next empty row).
$ws.Range("A8").Select()
